# Actualización automática 2025-10-15 12:30:08
# Inserts a new client row "CORONADO MONTERO LIDA VERONICA" (all zero values)
# right before "CORREA IGLESIAS RAMIRO MARCELO" in both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, shifting all following
# rows down by one, and updates the trailing summary rows accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A:R, data rows 2..344, summary row 344)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new blank row at row 286 (pushes CORREA IGLESIAS... and below down by 1)
$ws1.Rows.Item(286).Insert()

$ws1.Range("A286").Value = "OFICINA-CATAECSA"
$ws1.Range("B286").Value = "CORONADO MONTERO LIDA VERONICA"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(286, $c).Value = 0
}

# Fix up the trailing "X de 342" -> "X de 343" counts, now on row 345
$countsRow = 345
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item($countsRow, $c)
    $old = $cell.Value2
    $cell.Value = $old -replace "de 342", "de 343"
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A:G, data rows 2..348, summary row 348)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert a new blank row at row 290 (pushes CORREA IGLESIAS... and below down by 1)
$ws2.Rows.Item(290).Insert()

$ws2.Range("A290").Value = "OFICINA-CATAECSA"
$ws2.Range("B290").Value = "CORONADO MONTERO LIDA VERONICA"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(290, $c).Value = 0
}
